# Add the new "L6" sheet as the last tab (matches sheetId=7, after "Goal totals v2")
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "L6"

# Header row (B1:E1) - A1 stays blank
$ws.Range("B1").Value = "Form"
$ws.Range("C1").Value = "Goals scored"
$ws.Range("D1").Value = "Goals conceded"
$ws.Range("E1").Value = "Total Goals"

# Column B - Form (last 6 games), one row per team, alphabetical order
$ws.Range("B2").Value = "Aberdeen,L D L W D W"
$ws.Range("B3").Value = "Celtic,W D D W D L"
$ws.Range("B4").Value = "Dundee United,D D W W L L"
$ws.Range("B5").Value = "Hamilton,D L D L L W"
$ws.Range("B6").Value = "Hibernian,L W D L W L"
$ws.Range("B7").Value = "Kilmarnock,D L W D W L"
$ws.Range("B8").Value = "Livingston,L W D L L L"
$ws.Range("B9").Value = "Motherwell,W W L W W W"
$ws.Range("B10").Value = "Rangers,W W D W D W"
$ws.Range("B11").Value = "Ross County,W L L D L W"
$ws.Range("B12").Value = "St Johnstone,D W W L D W"
$ws.Range("B13").Value = "St Mirren,W L D L W L"

# Column C - Goals scored (last 6 games), one row per team, alphabetical order
$ws.Range("C2").Value = "Aberdeen,0 0 0 1 1 2"
$ws.Range("C3").Value = "Celtic,1 0 1 6 1 1"
$ws.Range("C4").Value = "Dundee United,1 0 1 1 0 0"
$ws.Range("C5").Value = "Hamilton,0 1 1 0 0 2"
$ws.Range("C6").Value = "Hibernian,0 2 1 1 2 0"
$ws.Range("C7").Value = "Kilmarnock,1 2 4 2 3 0"
$ws.Range("C8").Value = "Livingston,1 2 1 0 1 1"
$ws.Range("C9").Value = "Motherwell,2 3 1 1 1 2"
$ws.Range("C10").Value = "Rangers,1 3 1 2 1 4"
$ws.Range("C11").Value = "Ross County,3 1 0 2 1 2"
$ws.Range("C12").Value = "St Johnstone,1 1 1 0 1 1"
$ws.Range("C13").Value = "St Mirren,1 0 1 0 3 1"

# Column D - Goals conceded (last 6 games), one row per team, alphabetical order
$ws.Range("D2").Value = "Aberdeen,1 0 1 0 1 1"
$ws.Range("D3").Value = "Celtic,0 0 1 0 1 4"
$ws.Range("D4").Value = "Dundee United,1 0 0 0 3 2"
$ws.Range("D5").Value = "Hamilton,0 2 1 1 1 1"
$ws.Range("D6").Value = "Hibernian,1 1 1 2 1 1"
$ws.Range("D7").Value = "Kilmarnock,1 3 1 2 0 2"
$ws.Range("D8").Value = "Livingston,3 1 1 6 2 2"
$ws.Range("D9").Value = "Motherwell,0 1 4 0 0 0"
$ws.Range("D10").Value = "Rangers,0 0 1 1 1 1"
$ws.Range("D11").Value = "Ross County,2 2 1 2 3 0"
$ws.Range("D12").Value = "St Johnstone,1 0 0 1 1 0"
$ws.Range("D13").Value = "St Mirren,0 3 1 1 1 2"

# Column E - Total Goals (last 6 games), one row per team, alphabetical order
$ws.Range("E2").Value = "Aberdeen,1 0 1 1 2 3"
$ws.Range("E3").Value = "Celtic,1 0 2 6 2 5"
$ws.Range("E4").Value = "Dundee United,2 0 1 1 3 2"
$ws.Range("E5").Value = "Hamilton,0 3 2 1 1 3"
$ws.Range("E6").Value = "Hibernian,1 3 2 3 3 1"
$ws.Range("E7").Value = "Kilmarnock,2 5 5 4 3 2"
$ws.Range("E8").Value = "Livingston,4 3 2 6 3 3"
$ws.Range("E9").Value = "Motherwell,2 4 5 1 1 2"
$ws.Range("E10").Value = "Rangers,1 3 2 3 2 5"
$ws.Range("E11").Value = "Ross County,5 3 1 4 4 2"
$ws.Range("E12").Value = "St Johnstone,2 1 1 1 2 1"
$ws.Range("E13").Value = "St Mirren,1 3 2 1 4 3"

# Column A - league position numbers 1-12, written as text to match the
# convention used on the "Table" sheet (leading apostrophe forces text,
# then ClearFormats drops the quote-prefix style so the cell is plain text)
$ws.Range("A2").Value = "'1"
$ws.Range("A3").Value = "'2"
$ws.Range("A4").Value = "'3"
$ws.Range("A5").Value = "'4"
$ws.Range("A6").Value = "'5"
$ws.Range("A7").Value = "'6"
$ws.Range("A8").Value = "'7"
$ws.Range("A9").Value = "'8"
$ws.Range("A10").Value = "'9"
$ws.Range("A11").Value = "'10"
$ws.Range("A12").Value = "'11"
$ws.Range("A13").Value = "'12"
$ws.Range("A2:A13").ClearFormats()
